$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns keep their original plain-text representation (not numeric)
# by forcing Text number format before assigning the new values, matching the
# source data which stores prices/percentages as text strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.21%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.35%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.100"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.16%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07714"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.90%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.418"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.37%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.52%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "13.20%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1284"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.06%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1865"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.16%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09370"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.23%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04152"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.19%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.51%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001276"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.81%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005748"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.89%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.347"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.11%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3342"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.75%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.060"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.14%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1368"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.58%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.20%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04178"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.19%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001285"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.41%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004402"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "13.48%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001350"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.85%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02514"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.11%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05319"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.11%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005721"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-10.83%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007709"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.30%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.98%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007353"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.21%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007519"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-7.37%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3022"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.77%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.83%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04360"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-4.22%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.07%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
